$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 710.2308
$ws.Range("I19").Value = 428
$ws.Range("J19").Value = 886.625
$ws.Range("K19").Value = 428
$ws.Range("L19").Value = 886.625
$ws.Range("M19").Value = -253
$ws.Range("N19").Value = -1236.625
$ws.Range("H62").Value = 20834852
$ws.Range("I62").Value = 33334764
$ws.Range("J62").Value = 1665
$ws.Range("K62").Value = 33334764
$ws.Range("L62").Value = 1665
$ws.Range("M62").Value = -33334140
$ws.Range("N62").Value = -2913
$ws.Range("H65").Value = 20834852
$ws.Range("I65").Value = 33334764
$ws.Range("J65").Value = 1665
$ws.Range("K65").Value = 166673820
$ws.Range("L65").Value = 8325
$ws.Range("M65").Value = -166670700
$ws.Range("N65").Value = -14565
$ws.Range("H111").Value = 79184.62
$ws.Range("I111").Value = 2437.375
$ws.Range("J111").Value = 201980.2
$ws.Range("K111").Value = 7312.125
$ws.Range("L111").Value = 605940.6000000001
$ws.Range("M111").Value = -4245.125
$ws.Range("N111").Value = -612074.6000000001
$ws.Range("H112").Value = 1233.5
$ws.Range("I112").Value = 800
$ws.Range("J112").Value = 1305.75
$ws.Range("K112").Value = 2400
$ws.Range("L112").Value = 3917.25
$ws.Range("M112").Value = -1292
$ws.Range("N112").Value = -6133.25
$ws.Range("H113").Value = 20836550
$ws.Range("I113").Value = 3861
$ws.Range("J113").Value = 125000000
$ws.Range("K113").Value = 3861
$ws.Range("L113").Value = 125000000
$ws.Range("M113").Value = -607
$ws.Range("N113").Value = -125006508
$ws.Range("H115").Value = 1165.6666
$ws.Range("I115").Value = 748.5
$ws.Range("K115").Value = 2245.5
$ws.Range("M115").Value = -678.5
$ws.Range("H116").Value = 6235.84
$ws.Range("I116").Value = 8286
$ws.Range("J116").Value = 3160.6
$ws.Range("K116").Value = 8286
$ws.Range("L116").Value = 3160.6
$ws.Range("M116").Value = -4844
$ws.Range("N116").Value = -10044.6
$ws.Range("H118").Value = 996.0909
$ws.Range("I118").Value = 759.5
$ws.Range("K118").Value = 2278.5
$ws.Range("M118").Value = -621.5
$ws.Range("H137").Value = 1316.5319
$ws.Range("I137").Value = 1063.5278
$ws.Range("J137").Value = 2144.5454
$ws.Range("K137").Value = 3190.5834
$ws.Range("L137").Value = 6433.6362
$ws.Range("M137").Value = -640.5834000000004
$ws.Range("N137").Value = -11533.6362
$ws.Range("H141").Value = 1494.7885
$ws.Range("I141").Value = 1035.75
$ws.Range("J141").Value = 3024.9167
$ws.Range("K141").Value = 3107.25
$ws.Range("L141").Value = 9074.750100000001
$ws.Range("M141").Value = 2072.75
$ws.Range("N141").Value = -19434.7501

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1383.6364
$ws.Range("I2").Value = 1025.8572
$ws.Range("K2").Value = 1025.8572
$ws.Range("M2").Value = -912.8571999999999
$ws.Range("H32").Value = 5851.577
$ws.Range("I32").Value = 4095.8657
$ws.Range("J32").Value = 16545.455
$ws.Range("K32").Value = 4095.8657
$ws.Range("L32").Value = 16545.455
$ws.Range("M32").Value = -3808.8657
$ws.Range("N32").Value = -17119.455
$ws.Range("H45").Value = 15387.429
$ws.Range("I45").Value = 26103
$ws.Range("J45").Value = 1100
$ws.Range("K45").Value = 26103
$ws.Range("L45").Value = 1100
$ws.Range("M45").Value = -25726
$ws.Range("N45").Value = -1854
$ws.Range("H110").Value = 1723
$ws.Range("I110").Value = 1188.7778
$ws.Range("J110").Value = 2925
$ws.Range("K110").Value = 1188.7778
$ws.Range("L110").Value = 2925
$ws.Range("M110").Value = 856.2221999999999
$ws.Range("N110").Value = -7015
$ws.Range("H116").Value = 1383.6364
$ws.Range("I116").Value = 1025.8572
$ws.Range("K116").Value = 1025.8572
$ws.Range("M116").Value = 1268.1428

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1383.6364
$ws.Range("I3").Value = 1025.8572
$ws.Range("K3").Value = 1025.8572
$ws.Range("M3").Value = -911.8571999999999
$ws.Range("H99").Value = 142858850
$ws.Range("J99").Value = 2500
$ws.Range("L99").Value = 2500
$ws.Range("N99").Value = -5496
$ws.Range("H107").Value = 250001740
$ws.Range("J107").Value = 1000
$ws.Range("L107").Value = 1000
$ws.Range("N107").Value = -4840

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 5129893.5
$ws.Range("I16").Value = 7693690.5
$ws.Range("J16").Value = 2300
$ws.Range("K16").Value = 7693690.5
$ws.Range("L16").Value = 2300
$ws.Range("M16").Value = -7693403.5
$ws.Range("N16").Value = -2874
$ws.Range("H107").Value = 20834278
$ws.Range("I107").Value = 37037610
$ws.Range("K107").Value = 37037610
$ws.Range("M107").Value = -37035690
$ws.Range("H113").Value = 5129893.5
$ws.Range("I113").Value = 7693690.5
$ws.Range("J113").Value = 2300
$ws.Range("K113").Value = 7693690.5
$ws.Range("L113").Value = 2300
$ws.Range("M113").Value = -7691520.5
$ws.Range("N113").Value = -6640

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 98908.03
$ws.Range("I5").Value = 187.71794
$ws.Range("J5").Value = 273912.22
$ws.Range("K5").Value = 563.15382
$ws.Range("L5").Value = 821736.6599999999
$ws.Range("M5").Value = -451.15382
$ws.Range("N5").Value = -821960.6599999999
$ws.Range("H23").Value = 5000082.5
$ws.Range("I23").Value = 9090981
$ws.Range("J23").Value = 95.77778000000001
$ws.Range("K23").Value = 27272943
$ws.Range("L23").Value = 287.33334
$ws.Range("M23").Value = -27272708
$ws.Range("N23").Value = -757.33334
$ws.Range("H135").Value = 98908.03
$ws.Range("I135").Value = 187.71794
$ws.Range("J135").Value = 273912.22
$ws.Range("K135").Value = 1689.46146
$ws.Range("L135").Value = 2465209.98
$ws.Range("M135").Value = 845.53854
$ws.Range("N135").Value = -2470279.98

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 1000
$ws.Range("I5").Value = 1000
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 1000
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -888
$ws.Range("N5").ClearContents()
$ws.Range("H57").Value = 18250
$ws.Range("J57").Value = 18250
$ws.Range("L57").Value = 18250
$ws.Range("N57").Value = -19890
$ws.Range("H107").Value = 1906.8334
$ws.Range("I107").Value = 735.6667
$ws.Range("J107").Value = 3078
$ws.Range("K107").Value = 735.6667
$ws.Range("L107").Value = 3078
$ws.Range("M107").Value = 1184.3333
$ws.Range("N107").Value = -6918
$ws.Range("H113").Value = 125001780
$ws.Range("I113").Value = 333334620
$ws.Range("J113").Value = 2060
$ws.Range("K113").Value = 333334620
$ws.Range("L113").Value = 2060
$ws.Range("M113").Value = -333332450
$ws.Range("N113").Value = -6400
$ws.Range("H132").Value = 2247.4443
$ws.Range("I132").Value = 1926.7391
$ws.Range("J132").Value = 2814.8462
$ws.Range("K132").Value = 5780.2173
$ws.Range("L132").Value = 8444.5386
$ws.Range("M132").Value = -3250.2173
$ws.Range("N132").Value = -13504.5386

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H134").Value = 40000
$ws.Range("J134").Value = 40000
$ws.Range("L134").Value = 40000
$ws.Range("N134").Value = -50140

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 55167.168
$ws.Range("J2").Value = 55167.168
$ws.Range("L2").Value = 55167.168
$ws.Range("N2").Value = -55391.168
$ws.Range("H10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()
$ws.Range("H14").Value = 9231.154
$ws.Range("I14").Value = 7090.909
$ws.Range("J14").Value = 10800.667
$ws.Range("K14").Value = 7090.909
$ws.Range("L14").Value = 10800.667
$ws.Range("M14").Value = -6922.909
$ws.Range("N14").Value = -11136.667
$ws.Range("H100").Value = 524
$ws.Range("J100").Value = 650
$ws.Range("L100").Value = 1300
$ws.Range("N100").Value = -2382
$ws.Range("H107").Value = 61176904
$ws.Range("I107").Value = 90909490
$ws.Range("J107").Value = 6667166
$ws.Range("K107").Value = 272728470
$ws.Range("L107").Value = 20001498
$ws.Range("M107").Value = -272726550
$ws.Range("N107").Value = -20005338
$ws.Range("H132").Value = 1128.921
$ws.Range("I132").Value = 826.89655
$ws.Range("J132").Value = 2102.111
$ws.Range("K132").Value = 2480.68965
$ws.Range("L132").Value = 6306.333
$ws.Range("M132").Value = 49.31034999999974
$ws.Range("N132").Value = -11366.333

